# "Hoan thien Ngoai Tru" -- finalize outpatient test-data workbook:
#  - bump the record/test id and id-card numbers on both sheets
#  - update the shared insurance-card number text
#  - (re)apply bold + thin border formatting to the "Check" header row
#  - tidy up selections / the print margins on the "Check" sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Check")

# --- Data sheet -----------------------------------------------------------
$ws1.Range("A2").Value = 3000
$ws1.Range("E2").Value = 46200608000
$ws1.Range("X2").Value = "DN4127460130000"

# --- Check sheet ------------------------------------------------------------
$ws2.Range("A2").Value = 3000
$ws2.Range("C2").Value = "DN4127460130000"

# Re-apply the header formatting (bold font + thin box border, centered /
# top-aligned) on the "Check" sheet's header row.
$hdr = $ws2.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Reset the stale selection on "Check" back to the top-left cell.
$ws2.Range("A1").Select() | Out-Null

# Default print margins on the "Check" sheet (0.75"/0.75"/1"/1", 0.5" header/footer).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Restore "Data" as the active sheet/selection --------------------------
$ws1.Activate() | Out-Null
$ws1.Range("AB4").Select() | Out-Null
